# Apply the edits described by the commit "Upload excel files with prices"
# to the bread_coop_2022-07-11 worksheet.
#
# Summary of changes:
#  1. All data rows (2-396) get their timestamp (column P) updated from
#     "2022-07-11 18:30:42" to "2022-07-11 20:49:39".
#  2. Two products that went out of stock get " - Online kein Bestand"
#     inserted into their productAriaLabel (column N) text, rows 5 and 100.
#  3. A handful of rows were re-ordered in the source scrape: rows 241/242,
#     314/315 and 325/326 are pairwise swapped, and rows 331-334 are
#     rotated (331->334, 332->331, 333->332, 334->333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bulk-update the timestamp column for every data row -----------------
$ws.Range("P2:P396").Value = "2022-07-11 20:49:39"

# --- 2. Update the two "out of stock" product labels -------------------------
$ws.Range("N5").Value = "Ölz Super Soft Sandwich - Online kein Bestand 4.10 Schweizer Franken"
$ws.Range("N100").Value = "Betty Bossi Naturaplan Bio Kuchenteig rund ausgewallt Ø32cm - Online kein Bestand 1.95 Schweizer Franken"

# --- 3. Re-order a few rows that shifted position in the source scrape ------

# Swap rows 241 and 242
$row241 = $ws.Range("A241:P241").Value2
$row242 = $ws.Range("A242:P242").Value2
$ws.Range("A241:P241").Value = $row242
$ws.Range("A242:P242").Value = $row241

# Swap rows 314 and 315
$row314 = $ws.Range("A314:P314").Value2
$row315 = $ws.Range("A315:P315").Value2
$ws.Range("A314:P314").Value = $row315
$ws.Range("A315:P315").Value = $row314

# Swap rows 325 and 326
$row325 = $ws.Range("A325:P325").Value2
$row326 = $ws.Range("A326:P326").Value2
$ws.Range("A325:P325").Value = $row326
$ws.Range("A326:P326").Value = $row325

# Rotate rows 331-334: new331=old332, new332=old333, new333=old334, new334=old331
$row331 = $ws.Range("A331:P331").Value2
$row332 = $ws.Range("A332:P332").Value2
$row333 = $ws.Range("A333:P333").Value2
$row334 = $ws.Range("A334:P334").Value2
$ws.Range("A331:P331").Value = $row332
$ws.Range("A332:P332").Value = $row333
$ws.Range("A333:P333").Value = $row334
$ws.Range("A334:P334").Value = $row331

# Re-apply the new timestamp to the rows we just overwrote, since the
# Value2 snapshots above still carried the *old* timestamp string.
$ws.Range("P241:P242").Value = "2022-07-11 20:49:39"
$ws.Range("P314:P315").Value = "2022-07-11 20:49:39"
$ws.Range("P325:P326").Value = "2022-07-11 20:49:39"
$ws.Range("P331:P334").Value = "2022-07-11 20:49:39"
